$wb = $excel.ActiveWorkbook

# Add a new worksheet named "L6" after the last existing sheet ("Goal totals v2")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "L6"

# Header row (row 1): B1:E1
$ws.Range("B1").Value = "Form"
$ws.Range("C1").Value = "Goals scored"
$ws.Range("D1").Value = "Goals conceded"
$ws.Range("E1").Value = "Total Goals"

# Column A (row 2:25) holds the rank numbers 1-24 as TEXT (matches the
# "Table" sheet convention where A-column ranks are shared strings, not numbers).
# Pre-format as Text so Excel stores the typed digits as a string, not a number.
$ws.Range("A2:A25").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "1"
$ws.Cells.Item(3, 1).Value = "2"
$ws.Cells.Item(4, 1).Value = "3"
$ws.Cells.Item(5, 1).Value = "4"
$ws.Cells.Item(6, 1).Value = "5"
$ws.Cells.Item(7, 1).Value = "6"
$ws.Cells.Item(8, 1).Value = "7"
$ws.Cells.Item(9, 1).Value = "8"
$ws.Cells.Item(10, 1).Value = "9"
$ws.Cells.Item(11, 1).Value = "10"
$ws.Cells.Item(12, 1).Value = "11"
$ws.Cells.Item(13, 1).Value = "12"
$ws.Cells.Item(14, 1).Value = "13"
$ws.Cells.Item(15, 1).Value = "14"
$ws.Cells.Item(16, 1).Value = "15"
$ws.Cells.Item(17, 1).Value = "16"
$ws.Cells.Item(18, 1).Value = "17"
$ws.Cells.Item(19, 1).Value = "18"
$ws.Cells.Item(20, 1).Value = "19"
$ws.Cells.Item(21, 1).Value = "20"
$ws.Cells.Item(22, 1).Value = "21"
$ws.Cells.Item(23, 1).Value = "22"
$ws.Cells.Item(24, 1).Value = "23"
$ws.Cells.Item(25, 1).Value = "24"

# Data rows 2:25, columns B (Form), C (Goals scored), D (Goals conceded), E (Total Goals)
# Form -> column B
$ws.Range("B2").Value = "Accrington,D L W D D D"
$ws.Range("B3").Value = "AFC Wimbledon,W W W D D L"
$ws.Range("B4").Value = "Blackpool,D W L L W W"
$ws.Range("B5").Value = "Bristol Rvs,W D L L L L"
$ws.Range("B6").Value = "Burton,W D L D W D"
$ws.Range("B7").Value = "Charlton,W D W L D D"
$ws.Range("B8").Value = "Crewe,D L W D D W"
$ws.Range("B9").Value = "Doncaster,L W L L D W"
$ws.Range("B10").Value = "Fleetwood Town,W L L W L D"
$ws.Range("B11").Value = "Gillingham,L D L W D D"
$ws.Range("B12").Value = "Hull,W W W D W W"
$ws.Range("B13").Value = "Ipswich,D L D L D W"
$ws.Range("B14").Value = "Lincoln,W W W L W D"
$ws.Range("B15").Value = "Milton Keynes Dons,D L W W W D"
$ws.Range("B16").Value = "Northampton,L D L W D L"
$ws.Range("B17").Value = "Oxford,W W W L W W"
$ws.Range("B18").Value = "Peterboro,W W L W D D"
$ws.Range("B19").Value = "Plymouth,L L D L L L"
$ws.Range("B20").Value = "Portsmouth,D L L W D W"
$ws.Range("B21").Value = "Rochdale,W W W D D L"
$ws.Range("B22").Value = "Shrewsbury,L L L W L L"
$ws.Range("B23").Value = "Sunderland,L L D D L W"
$ws.Range("B24").Value = "Swindon,L L L W L L"
$ws.Range("B25").Value = "Wigan,W W W W D L"

# Goals scored -> column C
$ws.Range("C2").Value = "Accrington,0 1 2 3 3 1"
$ws.Range("C3").Value = "AFC Wimbledon,3 4 2 0 3 1"
$ws.Range("C4").Value = "Blackpool,0 1 0 0 1 3"
$ws.Range("C5").Value = "Bristol Rvs,2 1 0 0 0 0"
$ws.Range("C6").Value = "Burton,3 1 0 1 5 1"
$ws.Range("C7").Value = "Charlton,2 0 6 0 2 1"
$ws.Range("C8").Value = "Crewe,0 0 2 1 2 1"
$ws.Range("C9").Value = "Doncaster,0 2 1 0 2 2"
$ws.Range("C10").Value = "Fleetwood Town,1 1 0 1 2 1"
$ws.Range("C11").Value = "Gillingham,1 0 2 1 2 1"
$ws.Range("C12").Value = "Hull,3 3 2 2 2 3"
$ws.Range("C13").Value = "Ipswich,0 0 0 0 0 2"
$ws.Range("C14").Value = "Lincoln,4 1 1 1 1 3"
$ws.Range("C15").Value = "Milton Keynes Dons,0 0 1 2 5 1"
$ws.Range("C16").Value = "Northampton,0 1 1 3 2 0"
$ws.Range("C17").Value = "Oxford,6 4 3 1 3 3"
$ws.Range("C18").Value = "Peterboro,3 3 0 1 2 3"
$ws.Range("C19").Value = "Plymouth,0 0 1 0 1 1"
$ws.Range("C20").Value = "Portsmouth,0 0 1 1 3 3"
$ws.Range("C21").Value = "Rochdale,2 3 1 1 3 1"
$ws.Range("C22").Value = "Shrewsbury,1 0 1 1 0 2"
$ws.Range("C23").Value = "Sunderland,1 0 2 3 0 3"
$ws.Range("C24").Value = "Swindon,0 1 1 3 0 1"
$ws.Range("C25").Value = "Wigan,4 2 2 2 1 1"

# Goals conceded -> column D
$ws.Range("D2").Value = "Accrington,0 3 1 3 3 1"
$ws.Range("D3").Value = "AFC Wimbledon,0 1 1 0 3 3"
$ws.Range("D4").Value = "Blackpool,0 0 1 1 0 0"
$ws.Range("D5").Value = "Bristol Rvs,1 1 1 2 1 1"
$ws.Range("D6").Value = "Burton,0 1 1 1 2 1"
$ws.Range("D7").Value = "Charlton,1 0 0 1 2 1"
$ws.Range("D8").Value = "Crewe,0 2 0 1 2 0"
$ws.Range("D9").Value = "Doncaster,3 0 2 1 2 1"
$ws.Range("D10").Value = "Fleetwood Town,0 2 2 0 5 1"
$ws.Range("D11").Value = "Gillingham,4 0 3 0 2 1"
$ws.Range("D12").Value = "Hull,0 0 1 2 1 1"
$ws.Range("D13").Value = "Ipswich,0 3 0 3 0 1"
$ws.Range("D14").Value = "Lincoln,0 0 0 2 0 3"
$ws.Range("D15").Value = "Milton Keynes Dons,0 4 0 0 0 1"
$ws.Range("D16").Value = "Northampton,3 1 3 0 2 3"
$ws.Range("D17").Value = "Oxford,0 1 2 2 1 2"
$ws.Range("D18").Value = "Peterboro,0 1 1 0 2 3"
$ws.Range("D19").Value = "Plymouth,3 3 1 6 3 3"
$ws.Range("D20").Value = "Portsmouth,0 1 3 0 3 1"
$ws.Range("D21").Value = "Rochdale,1 1 0 1 3 2"
$ws.Range("D22").Value = "Shrewsbury,4 2 2 0 1 3"
$ws.Range("D23").Value = "Sunderland,2 1 2 3 1 1"
$ws.Range("D24").Value = "Swindon,3 2 4 1 5 2"
$ws.Range("D25").Value = "Wigan,1 1 0 1 1 3"

# Total Goals form -> column E
$ws.Range("E2").Value = "Accrington,0 4 3 6 6 2"
$ws.Range("E3").Value = "AFC Wimbledon,3 5 3 0 6 4"
$ws.Range("E4").Value = "Blackpool,0 1 1 1 1 3"
$ws.Range("E5").Value = "Bristol Rvs,3 2 1 2 1 1"
$ws.Range("E6").Value = "Burton,3 2 1 2 7 2"
$ws.Range("E7").Value = "Charlton,3 0 6 1 4 2"
$ws.Range("E8").Value = "Crewe,0 2 2 2 4 1"
$ws.Range("E9").Value = "Doncaster,3 2 3 1 4 3"
$ws.Range("E10").Value = "Fleetwood Town,1 3 2 1 7 2"
$ws.Range("E11").Value = "Gillingham,5 0 5 1 4 2"
$ws.Range("E12").Value = "Hull,3 3 3 4 3 4"
$ws.Range("E13").Value = "Ipswich,0 3 0 3 0 3"
$ws.Range("E14").Value = "Lincoln,4 1 1 3 1 6"
$ws.Range("E15").Value = "Milton Keynes Dons,0 4 1 2 5 2"
$ws.Range("E16").Value = "Northampton,3 2 4 3 4 3"
$ws.Range("E17").Value = "Oxford,6 5 5 3 4 5"
$ws.Range("E18").Value = "Peterboro,3 4 1 1 4 6"
$ws.Range("E19").Value = "Plymouth,3 3 2 6 4 4"
$ws.Range("E20").Value = "Portsmouth,0 1 4 1 6 4"
$ws.Range("E21").Value = "Rochdale,3 4 1 2 6 3"
$ws.Range("E22").Value = "Shrewsbury,5 2 3 1 1 5"
$ws.Range("E23").Value = "Sunderland,3 1 4 6 1 4"
$ws.Range("E24").Value = "Swindon,3 3 5 4 5 3"
$ws.Range("E25").Value = "Wigan,5 3 2 3 2 4"

# Restore the originally active sheet/tab selection
$wb.Worksheets.Item(1).Activate()
